$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns keep text formatting (avoid Excel auto-converting
# numeric-looking strings like "1.00" or "84.00" into numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.349.14'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.747.24'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.37'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.32'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.745.60'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.19'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.376.84'
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.750.53'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.409.18'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.87'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.64'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '466.25'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.00'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('E25').Value = '  +7.37%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.897.93'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.77'
$ws.Range('E31').Value = '  -3.73%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.29'
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.83'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.16'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.17'
$ws.Range('E35').Value = '  +1.56%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  --%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.702.18'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.101'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  -3.45%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.138'
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.78'
$ws.Range('E42').Value = '  +0.98%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.78'
$ws.Range('E45').Value = '  +14.40%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.301'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.55'
$ws.Range('E47').Value = '  +3.37%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.91'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.46'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '391.42'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.47'
$ws.Range('E51').Value = '  +0.35%  '
